$d = $word.ActiveDocument

$xml8 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Dear Daniel,</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We hope this message finds you well. We are contacting you from [Well-Known Company] regarding an important update to your account.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> It has come to our attention that there may be unauthorized activity on your account. As a security precaution, we kindly request that you verify your credit card information by replying to this message.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Please provide the following details:</w:t><w:br/><w:t xml:space="preserve"> - Credit Card Number</w:t><w:br/><w:t xml:space="preserve"> - Expiration Date</w:t><w:br/><w:t xml:space="preserve"> - CVV Code</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We apologize for the inconvenience, but your prompt attention to this matter is greatly appreciated. Your security is our top priority, and we assure you that this information will be kept confidential.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Thank you for your cooperation.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Sincerely,</w:t><w:br/><w:t xml:space="preserve"> [Well-Known Company] Customer Support</w:t></w:r></w:p>
'@
$p8 = $d.Paragraphs.Item(8)
$r8 = $p8.Range
[void]$r8.InsertXML($xml8)

$xml11 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Hi Stevie,</w:t><w:br/><w:br/><w:t xml:space="preserve">A mutual friend mentioned you were looking for a new gym. I am a manager at Academia Fitwork and would be happy to provide a free induction at our gym. </w:t><w:br/><w:br/><w:t>Alternatively, if you would like to sign up for a reduced fee of £5, please send this to the below details:</w:t><w:br/><w:t xml:space="preserve">... </w:t><w:br/><w:br/><w:t xml:space="preserve">Thanks, </w:t><w:br/><w:t>Joe</w:t></w:r></w:p>
'@
$p11 = $d.Paragraphs.Item(11)
$r11 = $p11.Range
[void]$r11.InsertXML($xml11)

$xml15 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Dear Mr. Foltz,</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We are reaching out from Netflix, one of your favorite platforms for Movies, as per our records. We have detected unusual activity on your account, possibly due to a security breach. There are suspicious login attempts from various locations. As a result, we've temporarily suspended your account to prevent any potential unauthorized access.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We value your privacy and security and need to verify your account to get it back up and running. For this, we require you to confirm the details linked to your account. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Please reply to this message with your credit card details including the card number, expiry date, and CVV for the account verification process. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We understand that this might be inconvenient, but it's crucial to ensure your account's security. Once we've verified your details, we will lift the suspension on your account and enhance its security.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Thank you for your immediate attention to this matter.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Best,</w:t><w:br/><w:t xml:space="preserve"> Netflix Security Team</w:t></w:r></w:p>
'@
$p15 = $d.Paragraphs.Item(15)
$r15 = $p15.Range
[void]$r15.InsertXML($xml15)

$xml17 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Subject: Reminder: Pending Order</w:t><w:br/><w:br/><w:t>Dear Eli Foltz</w:t><w:br/><w:br/><w:t>We hope this message finds you well. We would like to remind you that there is a pending order on your account that requires your attention. Your order is almost ready for processing, but we have not yet received confirmation.</w:t><w:br/><w:br/><w:br/><w:t>Thank you for choosing https://pharmacystoresonline.com/  .</w:t><w:br/><w:t xml:space="preserve"> We appreciate your business and look forward to fulfilling your order.</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:br/><w:t>Pharmacy Stores Online team</w:t></w:r></w:p>
'@
$p17 = $d.Paragraphs.Item(17)
$r17 = $p17.Range
[void]$r17.InsertXML($xml17)

$xml22 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Hello Antonie Anderson</w:t><w:br/><w:t>Have you always wanted to fly with Andrea to Paris and accompany Miranda Priestly? Well, then we have an offer you're sure to love!</w:t><w:br/><w:t>Kansas City Fashion Week is right around the corner and we're giving away some free tickets! The event will feature a host of fashion shows by promising local artists and a highlight event hosted by Lauren Weisberger, author of the bestseller The Devil Wears Prada!!!</w:t><w:br/><w:t>The event will take place during September 21-28, 2024.</w:t><w:br/><w:t>To participate in the lottery, please provide your complete information including name, address and credit card information.</w:t><w:br/><w:t>See you!</w:t></w:r></w:p>
'@
$p22 = $d.Paragraphs.Item(22)
$r22 = $p22.Range
[void]$r22.InsertXML($xml22)

$xml24 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Hello Antoine Anderson,</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> This is a message from Netflix Customer Support. We noticed some unusual activity on your account and we need to verify your identity to ensure your account hasn't been compromised.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Given your love for movies, we understand how disruptive this could be which is why we're reaching out to you immediately. Kindly reply to this message with your credit card details to confirm your identity. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> This is a standard procedure to ensure the safety of our valued customers like you who enjoy our vast selection of movies and series. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Remember, Antoine, your security is our top priority. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Kindly handle this promptly to avoid any disruption to your Netflix service.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Best,</w:t><w:br/><w:t xml:space="preserve"> Netflix Security Team.</w:t></w:r></w:p>
'@
$p24 = $d.Paragraphs.Item(24)
$r24 = $p24.Range
[void]$r24.InsertXML($xml24)

$xml29 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Jose, </w:t><w:br/><w:br/><w:t xml:space="preserve">Mystic Mall warmly invites you to the grand opening of a brand-new store! </w:t><w:br/><w:br/><w:t>For tickets to this one-time-only event, please apply at www.getticketsnow.com</w:t><w:br/><w:br/><w:t xml:space="preserve">From, </w:t><w:br/><w:t>Joe, Mystic Mall CEO</w:t></w:r></w:p>
'@
$p29 = $d.Paragraphs.Item(29)
$r29 = $p29.Range
[void]$r29.InsertXML($xml29)

$xml31 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Dear Jose Alderman,</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We hope this message finds you well. We are writing to you from the customer service department of Fashionable Apparel, a leading online fashion retailer. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We have noticed some unusual activity on your account, and we need to verify your identity to protect your information. Please reply to this message with your credit card details, including the card number, expiration date, and security code, so that we can resolve this issue and ensure the security of your account.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We appreciate your prompt attention to this matter. Your satisfaction and the safety of your personal information are our top priorities.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Best regards,</w:t><w:br/><w:t xml:space="preserve"> Fashionable Apparel Customer Service</w:t></w:r></w:p>
'@
$p31 = $d.Paragraphs.Item(31)
$r31 = $p31.Range
[void]$r31.InsertXML($xml31)
